# LQDA_YR_FIN.xlsx update: insert a new "current period" column (D) in front
# of the existing Income Statement / Balance Sheet / Cash Flow Statement
# blocks, pushing the previously-entered years one column to the right, and
# populate the new column with the freshly scraped figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a blank column before column D. This shifts D:K -> E:L for every
#    row, carrying values/styles/number-formats along with them (exactly
#    like pressing "Insert Sheet Columns" in the UI with D selected).
$ws.Columns("D:D").Insert()

# 2) The freshly inserted column D is blank and (because Excel clones the
#    format of the column to its left on insert) picked up column C's
#    style. Re-stamp it with the correct per-row style/number-format by
#    copying formats from the now-adjacent column E (which holds exactly
#    what used to be column D, style included).
$ws.Range("E7:E102").Copy() | Out-Null
$ws.Range("D7:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3) Drop in the new period's figures in column D for every data row that
#    carries a value (separator / spacer rows stay blank, exactly as they
#    were for the other periods).
$newValues = @{
    7  = 43465
    8  = 2700
    9  = 100
    10 = 2600
    12 = 28700
    13 = 0
    14 = -100
    15 = 0
    17 = 37400
    18 = -34700
    20 = 600
    21 = -32600
    22 = 19000
    23 = -53100
    24 = 0
    25 = 0
    26 = -53100
    27 = -53100
    28 = 0
    29 = 0
    30 = 0
    31 = 0
    32 = -600
    33 = -53100
    34 = 0
    35 = -53100
    38 = 43465
    41 = 39500
    42 = 0
    43 = 300
    44 = 0
    45 = 200
    46 = 40000
    47 = 0
    48 = 8100
    49 = 0
    50 = 0
    51 = 0
    52 = 1300
    53 = 0
    54 = 49400
    57 = 3200
    58 = 800
    59 = 4200
    60 = 8200
    61 = 12000
    62 = 10500
    63 = 0
    64 = 0
    65 = 0
    66 = 30700
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = -167100
    73 = 0
    74 = 0
    75 = 0
    76 = 18700
    77 = 0
    80 = 43465
    81 = -53100
    83 = 1500
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = -31800
    91 = -900
    92 = 0
    93 = 0
    94 = -900
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = 68800
    101 = 0
    102 = 36100
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 4).Value2 = $newValues[$row]
}

# 4) Row 14 ("Non Recurring") now only reports the new period's figure; the
#    previously-scraped periods that used to read 0 are re-marked "NA" just
#    like the rest of that row already was, except the oldest column (now
#    K, the 2013 figure) which still reports an actual 0.
foreach ($col in @("E", "F", "G", "H", "I", "J")) {
    $ws.Range($col + "14").Value2 = "NA"
}
$ws.Range("K14").Value2 = 0
